$d = $word.ActiveDocument

# --- Step 1: the stray "ß" run right after "Discretion?" is removed. ---
$null = $d.Content.Find.Execute("ß", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# --- Step 2: the hidden "_GoBack" bookmark currently sits at the end of the
#     "Discretion?" paragraph; it needs to move to the end of the new last
#     list item, so drop it here and re-create it below. ---
try {
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete()
} catch {
    # no existing _GoBack bookmark - nothing to remove
}

# --- Step 3: find the "Discretion?" list paragraph and the (now empty)
#     trailing list paragraph that immediately follows it. ---
$discIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($discIndex -eq -1 -and $p.Range.Text -like "Discretion?*") {
        $discIndex = $i
    }
}

$discPara = $d.Paragraphs($discIndex)
$trailingPara = $d.Paragraphs($discIndex + 1)

# Range spanning from right before the "Discretion?" paragraph mark through
# the end of the trailing empty list paragraph - this is everything that
# needs replacing with the three new FAQ entries.
$rStart = $discPara.Range.End - 1
$rEnd = $trailingPara.Range.End
$r = $d.Range($rStart, $rEnd)

$newItemsXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve"> Why so volatile?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve"> Why should I use your model? It looks hard&#8230;</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve"> Do you really believe the CB should respond so much with the interest rate?</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$r.InsertXML($newItemsXml)
